# Updates the cryptos price/volume snapshot (GitHub Actions scrape refresh).
# Price values that look like plain numbers are set with a leading
# apostrophe so Excel stores them as text (preserving exact formatting,
# e.g. trailing zeros like "32.00" or "1.00", instead of coercing to a
# numeric type). Values containing multiple "." (European-style grouping,
# e.g. "67.310.52") are never auto-numeric in Excel so need no prefix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.310.52"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "3.518.08"
$ws.Range("E3").Value = "  -0.86%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'610.22"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").Value = "'150.99"
$ws.Range("E6").Value = "  -1.95%  "

$ws.Range("D7").Value = "3.516.56"
$ws.Range("E7").Value = "  -0.79%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.481"
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("E10").Value = "  -0.92%  "

$ws.Range("D11").Value = "'7.06"
$ws.Range("E11").Value = "  +2.54%  "

$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = "  -1.38%  "

$ws.Range("D14").Value = "4.114.26"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").Value = "'32.00"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("D16").Value = "3.518.11"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "67.349.27"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("D20").Value = "'15.23"
$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("D21").Value = "'444.22"
$ws.Range("E21").Value = "  -2.03%  "

$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").Value = "'0.625"
$ws.Range("E23").Value = "  -2.27%  "

$ws.Range("D24").Value = "'77.35"
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("E25").Value = "  +9.07%  "

$ws.Range("D26").Value = "3.658.57"
$ws.Range("E26").Value = "  -0.95%  "

$ws.Range("D27").Value = "'10.34"
$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  +0.67%  "

$ws.Range("E30").Value = "  -2.52%  "

$ws.Range("D31").Value = "'1.54"
$ws.Range("E31").Value = "  -4.92%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "'0.163"
$ws.Range("E33").Value = "  +3.10%  "

$ws.Range("D34").Value = "'25.88"
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("D35").Value = "'6.17"
$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("D36").Value = "3.509.24"
$ws.Range("E36").Value = "  -1.18%  "

$ws.Range("E37").Value = "  -3.01%  "

$ws.Range("D38").Value = "'8.01"
$ws.Range("E38").Value = "  -0.62%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'177.68"
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  +4.89%  "

$ws.Range("D43").Value = "'0.0880"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("D44").Value = "'5.44"
$ws.Range("E44").Value = "  -3.07%  "

$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("D46").Value = "'28.79"
$ws.Range("E46").Value = "  -0.69%  "

$ws.Range("D47").Value = "'44.79"
$ws.Range("E47").Value = "  -2.33%  "

$ws.Range("D48").Value = "'2.61"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("D49").Value = "'1.26"
$ws.Range("E49").Value = "  +4.53%  "

$ws.Range("D50").Value = "'7.59"
$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("E51").Value = "  -1.03%  "
